$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '50.946.01'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '2.933.64'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '374.11'
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.63'
$ws.Range("E6").Value = '  -2.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.534'
$ws.Range("E7").Value = '  -1.42%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -1.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.40'
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0836'
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("D13").Value = '3.394.63'
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.90'
$ws.Range("E14").Value = '  -2.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.32'
$ws.Range("E15").Value = '  -1.70%  '
$ws.Range("D16").Value = '2.941.27'
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.973'
$ws.Range("E17").Value = '  +2.21%  '
$ws.Range("D18").Value = '50.959.73'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("E19").Value = '  -5.52%  '
$ws.Range("E20").Value = '  -2.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.52'
$ws.Range("E21").Value = '  -2.64%  '
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '263.94'
$ws.Range("E23").Value = '  +1.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.29'
$ws.Range("E24").Value = '  -0.86%  '
$ws.Range("E25").Value = '  +2.64%  '
$ws.Range("B26").Value = 'Filecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.15'
$ws.Range("E26").Value = '  +9.85%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.81'
$ws.Range("E27").Value = '  +9.30%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.167'
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.58'
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.83'
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.87'
$ws.Range("E33").Value = '  -0.51%  '
$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0450'
$ws.Range("E34").Value = '  +1.04%  '
$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '33.50'
$ws.Range("E35").Value = '  -2.89%  '
$ws.Range("E36").Value = '  -3.33%  '
$ws.Range("E37").Value = '  -0.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.98'
$ws.Range("E38").Value = '  -2.54%  '
$ws.Range("E39").Value = '  -0.90%  '
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("E41").Value = '  -4.13%  '
$ws.Range("E42").Value = '  -2.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '120.49'
$ws.Range("E43").Value = '  -2.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.288'
$ws.Range("E44").Value = '  +2.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.96'
$ws.Range("E45").Value = '  -4.95%  '
$ws.Range("E46").Value = '  -1.20%  '
$ws.Range("E47").Value = '  -2.38%  '
$ws.Range("E48").Value = '  +0.77%  '
$ws.Range("D49").Value = '1.970.77'
$ws.Range("E49").Value = '  -3.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0343'
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("E51").Value = '  -1.42%  '
